$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new columns I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Row 2 data
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

# Row 3 data
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 6
